{"js": "// \"Added a few more slots\" \u2014 move the meta-description blurb from the\n// top of the doc down to a caption under the feature-image prompt, and\n// relabel the heading that now precedes it.\n//\n// 1. Remove the whole \"Meta description: ...\" paragraph near the top.\n// 2. Insert a new bold \"Play Blazing Goddess for Free - Review of\n//    Lightning Box's Slot Game\" paragraph right before the final\n//    \"Prompt: ...\" paragraph.\n// 3. Replace the \"Prompt: ...\" paragraph's text with the tail of the old\n//    meta description (the part after \"Meta description: \"), keeping the\n//    paragraph's existing (italic) run formatting.\n\nconst body = context.document.body;\n\n// --- 1. Delete the \"Meta description\" paragraph -----------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst metaParagraph = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Meta description\") === 0\n);\nif (metaParagraph) {\n  metaParagraph.delete();\n  await context.sync();\n}\n\n// --- 2. Insert the new bold heading paragraph before \"Prompt: ...\" ----\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\n// Anchor on the paragraph right before the \"Prompt: ...\" paragraph (i.e.\n// two paragraphs before the end) and insert after it. Anchoring on the\n// paragraph *before* the prompt \u2014 rather than inserting \"before\" the\n// prompt paragraph itself \u2014 keeps the new paragraph from inheriting the\n// prompt paragraph's italic run formatting.\nconst anchorParagraph =\n  refreshedParagraphs.items[refreshedParagraphs.items.length - 2];\nconst anchorRange = anchorParagraph.getRange(Word.RangeLocation.end);\n\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>\" +\n  \"<w:t>Play Blazing Goddess for Free - Review of Lightning Box's Slot Game</w:t>\" +\n  \"</w:r></w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nanchorRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// --- 3. Replace the \"Prompt: ...\" paragraph's text ---------------------\nconst oldPromptText =\n  \"Prompt: Create a feature image for Blazing Goddess that reflects the \" +\n  \"game's exciting and adventurous vibe. The image should be in cartoon \" +\n  \"style and feature a happy Maya warrior with glasses. The warrior \" +\n  \"should be holding a flaming torch in one hand and a coconut cocktail \" +\n  \"in the other. The backdrop should be a tropical paradise with palm \" +\n  \"trees, sand, and sea. The Blazing Goddess symbol should be \" +\n  \"prominently displayed in the background, with flames erupting from \" +\n  \"it. The overall aesthetic should be bright, colorful, and energetic \" +\n  \"to capture the excitement of playing the Blazing Goddess slot \" +\n  \"machine.\";\nconst newPromptText =\n  \"Read our review of Blazing Goddess, a visually stunning online slot \" +\n  \"game from Lightning Box, or play for free with up to 60 free spins.\";\n\nconst matches = body.search(oldPromptText, { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(newPromptText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Added a few more slots\" - move the meta-description blurb from the\n# top of the doc down to a caption under the feature-image prompt, and\n# relabel the heading that now precedes it.\n#\n# 1. Remove the whole \"Meta description: ...\" paragraph near the top.\n# 2. Insert a new bold \"Play Blazing Goddess for Free - Review of\n#    Lightning Box's Slot Game\" paragraph right before the final\n#    \"Prompt: ...\" paragraph.\n# 3. Replace the \"Prompt: ...\" paragraph's text with the tail of the old\n#    meta description (the part after \"Meta description: \"), keeping the\n#    paragraph's existing (italic) run formatting.\n\n$d = $word.ActiveDocument\n\n# --- 1. Delete the \"Meta description\" paragraph ------------------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.IndexOf(\"Meta description\") -eq 0) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 2. Insert the new bold heading paragraph before \"Prompt: ...\" -----\n# Anchor on the paragraph right before the \"Prompt: ...\" paragraph (the\n# last paragraph in the document) and insert after it. Anchoring there -\n# rather than inserting \"before\" the prompt paragraph itself - keeps the\n# new paragraph from inheriting the prompt paragraph's italic run\n# formatting.\n$count = $d.Paragraphs.Count\n$anchorParagraph = $d.Paragraphs.Item($count - 1)\n$anchorParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($count)\n$newParagraphOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blazing Goddess for Free - Review of Lightning Box''s Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$newParagraph.Range.InsertXML($newParagraphOoxml)\n\n# --- 3. Replace the \"Prompt: ...\" paragraph's text ----------------------\n$oldPromptText = \"Prompt: Create a feature image for Blazing Goddess that reflects the game's exciting and adventurous vibe. The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a flaming torch in one hand and a coconut cocktail in the other. The backdrop should be a tropical paradise with palm trees, sand, and sea. The Blazing Goddess symbol should be prominently displayed in the background, with flames erupting from it. The overall aesthetic should be bright, colorful, and energetic to capture the excitement of playing the Blazing Goddess slot machine.\"\n$newPromptText = \"Read our review of Blazing Goddess, a visually stunning online slot game from Lightning Box, or play for free with up to 60 free spins.\"\n\n$promptParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$promptRange = $promptParagraph.Range\n$find = $promptRange.Find\n$find.Text = $oldPromptText\n$find.Replacement.Text = $newPromptText\n$find.Execute([ref]$oldPromptText, $false, $false, $false, $false, $false, $true, 1, $false, $newPromptText, 2)\n"}
